$d = $word.ActiveDocument

# 1. Fix typo: "Not what access controls..." -> "Note what access controls..."
$d.Content.Find.Execute("Not what access", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Note what access", 2) | Out-Null

# 2. Remove the now-obsolete "Add non-throwing swap." bullet point entirely
#    (whole paragraph, including its trailing paragraph mark).
$findRange = $d.Content
$found = $findRange.Find.Execute("Add non-throwing swap.", $true, $false, $false, $false, $false, `
                                  $true, 1, $false, "", 0)
if ($found) {
    $delRange = $d.Range($findRange.Start, $findRange.End + 1)
    $delRange.Delete() | Out-Null
}
